$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores every value as plain text in the source
# sheet (coinranking.com formats them with dotted thousands separators,
# fixed decimal places, etc.), even when a value happens to look like a
# plain number (e.g. "131.90", "504.25"). Force those cells to Text
# format before assigning so Excel does not silently coerce them to
# Number and strip significant trailing zeros / reformat them.

$ws.Range('D2').Value = '56.777.48'
$ws.Range('E2').Value = '  +0.82%  '
$ws.Range('D3').Value = '2.389.93'
$ws.Range('E3').Value = '  +1.32%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '504.25'
$ws.Range('E5').Value = '  -1.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.90'
$ws.Range('E6').Value = '  +2.73%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = '2.399.31'
$ws.Range('E9').Value = '  +1.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0970'
$ws.Range('E10').Value = '  +1.67%  '
$ws.Range('E11').Value = '  -1.12%  '
$ws.Range('E12').Value = '  +1.80%  '
$ws.Range('E13').Value = '  -4.54%  '
$ws.Range('D14').Value = '2.816.73'
$ws.Range('E14').Value = '  +1.35%  '
$ws.Range('D15').Value = '56.720.54'
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.71'
$ws.Range('E16').Value = '  +1.29%  '
$ws.Range('E17').Value = '  +2.01%  '
$ws.Range('D18').Value = '2.388.17'
$ws.Range('E18').Value = '  +1.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.19'
$ws.Range('E19').Value = '  -0.88%  '
$ws.Range('E20').Value = '  -0.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '308.92'
$ws.Range('E21').Value = '  -0.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.26'
$ws.Range('E22').Value = '  +0.61%  '
$ws.Range('E23').Value = '  +0.52%  '
$ws.Range('E24').Value = '  +0.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.23'
$ws.Range('E25').Value = '  +0.52%  '
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('E27').Value = '  -4.10%  '
$ws.Range('E28').Value = '  -0.59%  '
$ws.Range('E29').Value = '  +2.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '175.04'
$ws.Range('D31').Value = '0.0₃0720'
$ws.Range('E31').Value = '  +1.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.67'
$ws.Range('E32').Value = '  -0.50%  '
$ws.Range('E33').Value = '  +0.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.87'
$ws.Range('E34').Value = '  -4.81%  '
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('E37').Value = '  +1.43%  '
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('E39').Value = '  +2.89%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.81'
$ws.Range('E40').Value = '  +3.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.817'
$ws.Range('E41').Value = '  +3.49%  '
$ws.Range('E42').Value = '  +1.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '133.13'
$ws.Range('E43').Value = '  +5.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.35'
$ws.Range('E44').Value = '  +0.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.81'
$ws.Range('E45').Value = '  +0.71%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '250.13'
$ws.Range('E47').Value = '  -1.83%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0907'
$ws.Range('E48').Value = '  +0.78%  '
$ws.Range('E50').Value = '  +1.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.07'
$ws.Range('E51').Value = '  +8.74%  '
